# Final Updated Data Done
# Shorten the hard-coded "C:\Users\gaura\Desktop\..." file paths on the
# "Settings" sheet down to the relative "FFR Data\..." paths, and move the
# active selection to B21 (top of the visible Settings list).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

$ws.Range("B14").Value = 'FFR Data\Added Data\OM13_Faults_Data.xlsx'
$ws.Range("B15").Value = 'FFR Data\Intermediate Data\Faults_Data.xlsx'
$ws.Range("B16").Value = 'FFR Data\Added Data\FFR Data Actual Part with Part code.xlsx'
$ws.Range("B17").Value = 'FFR Data\Exception Handling Folder\Faults_Data.xlsx'

$ws.Range("B21").Value = 'FFR Data\Added Data\AC - IW POP OM1_21 TO OM13_21.xlsx'
$ws.Range("B22").Value = 'FFR Data\Intermediate Data\Population Data.xlsx'
$ws.Range("B23").Value = 'FFR Data\Exception Handling Folder\Population Data.xlsx'

$ws.Range("B27").Value = 'FFR Data\Added Data\Master  In Warranty Population Vs  Complaints OM12 2022.xlsx'
$ws.Range("B28").Value = 'FFR Data\Intermediate Data\AC FFR Models.xlsx'
$ws.Range("B29").Value = 'FFR Data\Exception Handling Folder\AC FFR Models.xlsx'

$ws.Range("B33").Value = 'FFR Data\Added Data\FFR_Data.xlsx'
$ws.Range("B34").Value = 'FFR Data\Intermediate Data\FFR_Data.xlsx'

$ws.Range("B37").Value = 'FFR Data\Intermediate Data\AC FFR PPT.pptx'
$ws.Range("B38").Value = 'FFR Data\Exception Handling Folder\AC FFR OM08.pptx'

$ws.Range("B40").Value = 'FFR Data\Final Data'

$ws.Range("B44").Value = 'FFR Data\FFR VB Script\AddingLineGraph.vbs'
$ws.Range("B45").Value = 'FFR Data\FFR VB Script\AutofitColumn.vbs'
$ws.Range("B46").Value = 'FFR Data\FFR VB Script\DeleteSheets.vbs'
$ws.Range("B47").Value = 'FFR Data\FFR VB Script\Graph.vbs'
$ws.Range("B48").Value = 'FFR Data\FFR VB Script\Merge.vbs'
$ws.Range("B49").Value = 'FFR Data\FFR VB Script\NewSheet.vbs'
$ws.Range("B50").Value = 'FFR Data\FFR VB Script\NewTableGraphBorder.vbs'
$ws.Range("B51").Value = 'FFR Data\FFR VB Script\OneTableGraphBorder.vbs'
$ws.Range("B52").Value = 'FFR Data\FFR VB Script\ThreeTableGraphBorder.vbs'
$ws.Range("B53").Value = 'FFR Data\FFR VB Script\TwoTableGraphBorder.vbs'
$ws.Range("B54").Value = 'FFR Data\FFR VB Script\Unmerge.vbs'
$ws.Range("B55").Value = 'FFR Data\FFR VB Script\UpdateChart.vbs'
$ws.Range("B56").Value = 'FFR Data\FFR VB Script\ZoomOut.vbs'

# Move the visible selection/cursor (matches the updated sheetView in the diff).
$ws.Range("B21").Select()
